# Recalculo de los resultados, cambios de LSPMW
# Updates the ranking tables on the General, ARMA, ARIMA and SETAR sheets.

$wb = $excel.ActiveWorkbook

function Set-RankingRow {
    param(
        $ws,
        [int]$row,
        [string]$modelo,
        $victorias,
        $empates,
        $derrotas,
        $score,
        $pctVictorias
    )
    $ws.Cells.Item($row, 2).Value = $modelo
    $ws.Cells.Item($row, 3).Value = $victorias
    $ws.Cells.Item($row, 4).Value = $empates
    $ws.Cells.Item($row, 5).Value = $derrotas
    $ws.Cells.Item($row, 6).Value = $score
    $ws.Cells.Item($row, 7).Value = $pctVictorias
}

# ----- General -----
$ws = $wb.Worksheets.Item("General")
Set-RankingRow $ws 4 "MCPS"       4 2 2  2  50
Set-RankingRow $ws 5 "LSPMW"      3 3 2  1  37.5
Set-RankingRow $ws 6 "AV-MCPS"    3 3 2  1  37.5
Set-RankingRow $ws 7 "DeepAR"     2 3 3  -1 25
Set-RankingRow $ws 8 "EnCQR-LSTM" 2 1 5  -3 25

# ----- ARMA -----
$ws = $wb.Worksheets.Item("ARMA")
Set-RankingRow $ws 4  "MCPS"                4 2 2 2  50
Set-RankingRow $ws 5  "AV-MCPS"             4 2 2 2  50
Set-RankingRow $ws 6  "LSPM"                3 3 2 1  37.5
Set-RankingRow $ws 7  "Block Bootstrapping" 3 0 5 -2 37.5
Set-RankingRow $ws 8  "AREPD"               2 0 6 -4 25
Set-RankingRow $ws 9  "EnCQR-LSTM"          1 1 6 -5 12.5
Set-RankingRow $ws 10 "LSPMW"               0 0 8 -8 0

# ----- ARIMA -----
$ws = $wb.Worksheets.Item("ARIMA")
Set-RankingRow $ws 3 "LSPM"       7 0 1  6  87.5
Set-RankingRow $ws 4 "MCPS"       4 2 2  2  50
Set-RankingRow $ws 5 "LSPMW"      3 3 2  1  37.5
Set-RankingRow $ws 6 "AV-MCPS"    3 3 2  1  37.5
Set-RankingRow $ws 7 "DeepAR"     2 3 3  -1 25
Set-RankingRow $ws 8 "EnCQR-LSTM" 2 1 5  -3 25

# ----- SETAR -----
$ws = $wb.Worksheets.Item("SETAR")
Set-RankingRow $ws 2  "Block Bootstrapping" 5 3 0 5  62.5
Set-RankingRow $ws 3  "DeepAR"              5 3 0 5  62.5
Set-RankingRow $ws 4  "Sieve Bootstrap"     4 4 0 4  50
Set-RankingRow $ws 5  "AV-MCPS"             0 8 0 0  0
Set-RankingRow $ws 6  "LSPM"                1 5 2 -1 12.5
Set-RankingRow $ws 7  "LSPMW"               1 4 3 -2 12.5
Set-RankingRow $ws 8  "EnCQR-LSTM"          0 5 3 -3 0
Set-RankingRow $ws 9  "MCPS"                0 5 3 -3 0
Set-RankingRow $ws 10 "AREPD"               0 3 5 -5 0
